# "Primer commit de aplicacion de stock"
# Update the current stock ("Stock_actual") figure for the first product
# row (Shampoo Coco, row 4) from 10 down to 7 units, leaving the edited
# cell selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7

$ws.Range("B4").Select() | Out-Null
